$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Mark all test cases as not-running ("N") except TestCase_A1 (row 2),
# which stays "Y", and TestCase_A13 (row 14), which was already "N".
$ws.Range("C3:C13").Value = "N"
$ws.Range("C15:C16").Value = "N"

# Update the selection to reflect the range that was just edited.
$ws.Range("C3:C16").Select()
